$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("energy")

# Insert a new blank row above row 13 (shifts old rows 13:45 down to 14:46)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 (code "17" -> "Traditional biomass"), matching the
# styling pattern used by similar rows (e.g. row 8: A blank/s1, B text/s1, C text/s4)
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = "'17"
$ws.Range("C13").Value = "Traditional biomass"

# Copy the formatting (styles) from row 8, which has the same A/B/C pattern
$ws.Range("A8").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# Update the active selection / view to match the committed state
$ws.Range("B14").Select()

Write-Host "done"
